# Auto-generated edits applying the RPA datasets push 2023-09-12 commit.
# Updates IPO tracking table rows 2-21 on Sheet1 to the latest snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '비아이매트릭스'
$ws.Range("B2").Value = '2023.10.19~10.25'
$ws.Range("C2").Value = '9,100~11,000'
$ws.Range("E2").Value = 10920
$ws.Range("F2").Value = 'IBK투자증권'
$ws.Range("A3").Value = '유투바이오'
$ws.Range("B3").Value = '2023.10.18~10.19'
$ws.Range("C3").Value = '3,300~3,900'
$ws.Range("E3").Value = 3724
$ws.Range("F3").Value = '신한투자증권'
$ws.Range("A4").Value = '퀄리타스반도체'
$ws.Range("B4").Value = '2023.10.06~10.13'
$ws.Range("C4").Value = '13,000~15,000'
$ws.Range("E4").Value = 23400
$ws.Range("F4").Value = '한국투자증권'
$ws.Range("A5").Value = '컨텍'
$ws.Range("B5").Value = '2023.10.06~10.13'
$ws.Range("C5").Value = '20,300~22,500'
$ws.Range("E5").Value = 41818
$ws.Range("F5").Value = '대신증권'
$ws.Range("A6").Value = '신성에스티'
$ws.Range("B6").Value = '2023.09.22~10.04'
$ws.Range("C6").Value = '22,000~25,000'
$ws.Range("E6").Value = 44000
$ws.Range("F6").Value = '미래에셋증권'
$ws.Range("A7").Value = '퓨릿(구.신디프)'
$ws.Range("B7").Value = '2023.09.20~09.26'
$ws.Range("C7").Value = '8,800~10,700'
$ws.Range("E7").Value = 36405
$ws.Range("F7").Value = '미래에셋증권'
$ws.Range("A8").Value = '에이치엠씨아이비스팩6호'
$ws.Range("B8").Value = '2023.09.19~09.20'
$ws.Range("E8").Value = 8000
$ws.Range("F8").Value = '현대차증권'
$ws.Range("A9").Value = '에스엘에스바이오'
$ws.Range("B9").Value = '2023.09.18~09.22'
$ws.Range("C9").Value = '8,200~9,400'
$ws.Range("E9").Value = 6314
$ws.Range("F9").Value = '하나증권'
$ws.Range("A10").Value = '신한스팩11호'
$ws.Range("B10").Value = '2023.09.14~09.15'
$ws.Range("C10").Value = '2,000~2,000'
$ws.Range("E10").Value = 36000
$ws.Range("F10").Value = '신한투자증권'
$ws.Range("A11").Value = '레뷰코퍼레이션'
$ws.Range("B11").Value = '2023.09.11~09.15'
$ws.Range("C11").Value = '11,500~13,200'
$ws.Range("E11").Value = 25760
$ws.Range("F11").Value = '삼성증권'
$ws.Range("A12").Value = '두산로보틱스'
$ws.Range("B12").Value = '2023.09.11~09.15'
$ws.Range("C12").Value = '21,000~26,000'
$ws.Range("E12").Value = 340200
$ws.Range("F12").Value = '한국투자증권,미래에셋증권,NH투자증권,KB증권,크레디트스위스증권'
$ws.Range("A13").Value = '한싹'
$ws.Range("B13").Value = '2023.09.08~09.14'
$ws.Range("C13").Value = '8,900~11,000'
$ws.Range("E13").Value = 13350
$ws.Range("F13").Value = '케이비증권'
$ws.Range("A14").Value = '밀리의서재'
$ws.Range("B14").Value = '2023.09.07~09.13'
$ws.Range("C14").Value = '20,000~23,000'
$ws.Range("E14").Value = 30000
$ws.Range("F14").Value = '미래에셋증권'
$ws.Range("A15").Value = '인스웨이브시스템즈'
$ws.Range("B15").Value = '2023.09.06~09.12'
$ws.Range("C15").Value = '20,000~24,000'
$ws.Range("D15").Value = '-'
$ws.Range("E15").Value = 22000
$ws.Range("F15").Value = '신영증권'
$ws.Range("A16").Value = '아이엠티'
$ws.Range("B16").Value = '2023.09.06~09.12'
$ws.Range("C16").Value = '10,500~12,000'
$ws.Range("D16").Value = '-'
$ws.Range("E16").Value = 16590
$ws.Range("F16").Value = '유안타증권,유진투자증권'
$ws.Range("A17").Value = '상상인스팩4호'
$ws.Range("B17").Value = '2023.08.28~08.29'
$ws.Range("E17").Value = 9000
$ws.Range("F17").Value = '상상인증권'
$ws.Range("A18").Value = '한화플러스스팩4호'
$ws.Range("B18").Value = '2023.08.24~08.25'
$ws.Range("E18").Value = 9500
$ws.Range("F18").Value = '한화투자증권'
$ws.Range("A19").Value = '대신밸런스스팩16호'
$ws.Range("B19").Value = '2023.08.17~08.18'
$ws.Range("E19").Value = 13000
$ws.Range("F19").Value = '대신증권'
$ws.Range("A20").Value = '한국스팩12호'
$ws.Range("B20").Value = '2023.08.16~08.17'
$ws.Range("E20").Value = 8000
$ws.Range("F20").Value = '한국투자증권'
$ws.Range("A21").Value = '유안타스팩11호'
$ws.Range("B21").Value = '2023.08.16~08.17'
$ws.Range("C21").Value = '2,000~2,000'
$ws.Range("E21").Value = 10000
$ws.Range("F21").Value = '유안타증권'

# D21 target value '2000' is numeric-looking text (must stay a text cell, like
# the existing D17:D20 cells). Copy an existing text-formatted '2000' cell over
# instead of assigning .Value directly, so it is not auto-converted to a number.
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D21").PasteSpecial()
$excel.CutCopyMode = $false
